# Updated cryptos list values (price & 1h volume change) and reorders
# the MXToken / Aave rows, per the Oct 12 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Cells hold plain text (prices use
# dotted grouping like "1.564.45", and percentages keep their padding
# spaces), so force text format before assigning to avoid Excel
# reinterpreting them as numbers.
$updates = [ordered]@{
    'D2' = '26.899.45'
    'E2' = '  -0.94%  '
    'D3' = '1.566.08'
    'E3' = '  +0.64%  '
    'E4' = '  -0.24%  '
    'D5' = '206.14'
    'E5' = '  +0.03%  '
    'D6' = '0.486'
    'E6' = '  -0.27%  '
    'E7' = '  -0.24%  '
    'D8' = '21.84'
    'E8' = '  -1.11%  '
    'E9' = '  +0.03%  '
    'D10' = '0.0586'
    'E10' = '  -0.83%  '
    'D11' = '0.0865'
    'E11' = '  +0.63%  '
    'D12' = '1.787.41'
    'E12' = '  +0.47%  '
    'D13' = '1.568.55'
    'E13' = '  +0.71%  '
    'E14' = '  -0.57%  '
    'E15' = '  +0.15%  '
    'D16' = '26.887.13'
    'E16' = '  -0.97%  '
    'D17' = '61.36'
    'E17' = '  -2.24%  '
    'D18' = '215.76'
    'E18' = '  +0.74%  '
    'D19' = '7.39'
    'E19' = '  +2.44%  '
    'D20' = '0.0₃0683'
    'E20' = '  -0.01%  '
    'E21' = '  -0.18%  '
    'D22' = '4.14'
    'E22' = '  +1.13%  '
    'D23' = '9.21'
    'E23' = '  -1.26%  '
    'E24' = '  +0.89%  '
    'D25' = '153.62'
    'E25' = '  +1.15%  '
    'D26' = '6.72'
    'E26' = '  +2.23%  '
    'E27' = '  +0.79%  '
    'E28' = '  -0.21%  '
    'E29' = '  -0.42%  '
    'D30' = '0.0467'
    'E30' = '  +1.22%  '
    'E31' = '  -3.24%  '
    'E32' = '  -0.14%  '
    'D33' = '1.401.61'
    'E33' = '  +1.95%  '
    'E34' = '  +0.29%  '
    'E35' = '  -0.33%  '
    'E36' = '  -0.37%  '
    'D37' = '0.918'
    'E37' = '  -3.23%  '
    'E38' = '  -0.22%  '
    'D39' = '0.530'
    'E39' = '  +3.30%  '
    'D40' = '0.814'
    'E40' = '  +0.76%  '
    'E41' = '  -0.20%  '
    'D42' = '0.993'
    'E42' = '  +0.53%  '
    'D43' = '5.56'
    'E43' = '  +6.75%  '
    'E44' = '  +0.61%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '63.91'
    'E45' = '  +1.26%  '
    'B46' = 'MXToken'
    'C46' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D46' = '2.18'
    'E46' = '  +0.99%  '
    'D47' = '1.700.85'
    'E47' = '  +0.56%  '
    'D48' = '86.82'
    'E48' = '  +2.03%  '
    'D49' = '0.0502'
    'E49' = '  +2.19%  '
    'D50' = '0.0₇0973'
    'E50' = '  -1.38%  '
    'D51' = '0.0953'
    'E51' = '  +1.54%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
}
